$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 22.02.2022 06:15"

# D4: was an inline string "+0.8", becomes a true number 0.8
$ws.Range("D4").Value = 0.8

# E4: was an inline string "2022-02-22 06:00:09", becomes a numeric date
# serial (matching the style/number format already used by the other
# rows in column E, e.g. E2/E3 which carry style index 2 / format 165)
$ws.Range("E4").Value = 44614.25010416667
$ws.Range("E4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
